{"js": "const newValues = [\"36+59=\", \"90-54=\", \"99-2=\", \"66-7=\", \"69-14=\", \"59+35=\", \"35+56=\", \"60+11=\", \"23+11=\", \"44-22=\", \"55+39=\", \"10+78=\", \"76+6=\", \"59-11=\", \"59-5=\", \"28+14=\", \"3+58=\", \"42+2=\", \"9+42=\", \"80-58=\", \"55+6=\", \"25+18=\", \"7+85=\", \"15+55=\", \"7+33=\", \"53-3=\", \"68-65=\", \"90-53=\", \"76-49=\", \"45-43=\", \"31-0=\", \"70-10=\", \"43+45=\", \"48-41=\", \"73-58=\", \"82-69=\", \"31-6=\", \"8-5=\", \"83-35=\", \"47+17=\", \"68-36=\", \"42+55=\", \"12+78=\", \"17+13=\", \"73-21=\", \"16+79=\", \"61+22=\", \"47-40=\", \"39+37=\", \"93-21=\", \"87-27=\", \"84+5=\", \"89-11=\", \"67+28=\", \"7+37=\", \"69-50=\", \"89-65=\", \"16+73=\", \"80-47=\", \"6+91=\", \"54+36=\", \"43+21=\", \"37-22=\", \"65+27=\", \"47+12=\", \"79-3=\", \"94-1=\", \"18+22=\", \"46-35=\", \"22-20=\", \"14+36=\", \"35+27=\", \"86+4=\", \"28-6=\", \"71-22=\", \"39+15=\", \"14+75=\", \"38+27=\", \"61+26=\", \"40+58=\", \"65+8=\", \"46+15=\", \"41+26=\", \"56-47=\", \"74+1=\", \"99-53=\", \"81-71=\", \"95-8=\", \"57-49=\", \"33+54=\", \"54+17=\", \"27-2=\", \"3+62=\", \"40-10=\", \"44-24=\", \"64+6=\", \"83+4=\", \"0+92=\", \"31+29=\", \"52-23=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load cells for every row up front.\nfor (const row of rows.items) {\n  row.load(\"cells/items\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    const para = paragraphs.items[0];\n    para.insertText(newValues[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "$newValues = @(\"36+59=\",\"90-54=\",\"99-2=\",\"66-7=\",\"69-14=\",\"59+35=\",\"35+56=\",\"60+11=\",\"23+11=\",\"44-22=\",\"55+39=\",\"10+78=\",\"76+6=\",\"59-11=\",\"59-5=\",\"28+14=\",\"3+58=\",\"42+2=\",\"9+42=\",\"80-58=\",\"55+6=\",\"25+18=\",\"7+85=\",\"15+55=\",\"7+33=\",\"53-3=\",\"68-65=\",\"90-53=\",\"76-49=\",\"45-43=\",\"31-0=\",\"70-10=\",\"43+45=\",\"48-41=\",\"73-58=\",\"82-69=\",\"31-6=\",\"8-5=\",\"83-35=\",\"47+17=\",\"68-36=\",\"42+55=\",\"12+78=\",\"17+13=\",\"73-21=\",\"16+79=\",\"61+22=\",\"47-40=\",\"39+37=\",\"93-21=\",\"87-27=\",\"84+5=\",\"89-11=\",\"67+28=\",\"7+37=\",\"69-50=\",\"89-65=\",\"16+73=\",\"80-47=\",\"6+91=\",\"54+36=\",\"43+21=\",\"37-22=\",\"65+27=\",\"47+12=\",\"79-3=\",\"94-1=\",\"18+22=\",\"46-35=\",\"22-20=\",\"14+36=\",\"35+27=\",\"86+4=\",\"28-6=\",\"71-22=\",\"39+15=\",\"14+75=\",\"38+27=\",\"61+26=\",\"40+58=\",\"65+8=\",\"46+15=\",\"41+26=\",\"56-47=\",\"74+1=\",\"99-53=\",\"81-71=\",\"95-8=\",\"57-49=\",\"33+54=\",\"54+17=\",\"27-2=\",\"3+62=\",\"40-10=\",\"44-24=\",\"64+6=\",\"83+4=\",\"0+92=\",\"31+29=\",\"52-23=\")\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $tbl.Cell($r, $c)\n    $rng = $cell.Range\n    $rng.End = $rng.End - 1\n    $rng.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
